$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# ---------------------------------------------------------------------
# 1) Title shape (id=1 / "Titre 1"): "A Proposal for draft-ietf-nmop-terminology"
#    -> "Next Step for draft-ietf-nmop-terminology"
# ---------------------------------------------------------------------
$titleShp = $s.Shapes.Item(1)
$titleTr = $titleShp.TextFrame.TextRange
$rA = $titleTr.Characters(1, 2)
$rA.Text = "Next "
$rProposal = $titleTr.Characters(6, 8)
$rProposal.Text = "Step"

# ---------------------------------------------------------------------
# Content placeholder shape (id=2 / "Espace reserve du contenu 2")
# Edit paragraphs bottom-up so earlier (not-yet-touched) paragraph
# indices never shift under us.
# ---------------------------------------------------------------------
$bodyShp = $s.Shapes.Item(2)
$tr = $bodyShp.TextFrame.TextRange

# --- paragraph 7: "NMOP documents / can define context-specific terms, / if needed"
#     -> "NMOP documents can define " + bold-italic "context-specific terms, if needed"
$para7 = $tr.Paragraphs(7, 1)
$p7start = $para7.Start
$para7.Text = "TEMP_PLACEHOLDER_7"
$para7b = $tr.Paragraphs(7, 1)
$para7b.Text = "NMOP documents can define context-specific terms, if needed"
$plain7 = "NMOP documents can define "
$bi7 = "context-specific terms, if needed"
$sub7 = $tr.Characters($p7start + $plain7.Length, $bi7.Length)
$sub7.Font.Bold = 1
$sub7.Font.Italic = 1

# --- paragraph 6: "If no major comment is raised, run a WGLC by end of Jan 2025"
#     -> plain "If no major comment is raised, " + bold-italic "run a WGLC by end of Jan 2025"
$para6 = $tr.Paragraphs(6, 1)
$p6start = $para6.Start
$plain6 = "If no major comment is raised, "
$bi6 = "run a WGLC by end of Jan 2025"
$sub6 = $tr.Characters($p6start + $plain6.Length, $bi6.Length)
$sub6.Font.Bold = 1
$sub6.Font.Italic = 1

# --- paragraph 5: "Let the document open for two/three months"
#     -> plain "Let the document " + bold-italic "open for two/three months"
$para5 = $tr.Paragraphs(5, 1)
$p5start = $para5.Start
$plain5 = "Let the document "
$bi5 = "open for two/three months"
$sub5 = $tr.Characters($p5start + $plain5.Length, $bi5.Length)
$sub5.Font.Bold = 1
$sub5.Font.Italic = 1

# --- paragraph 3: "A Proposal for discussion" -> "A proposal for discussion"
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "TEMP_PLACEHOLDER_3"
$para3b = $tr.Paragraphs(3, 1)
$para3b.Text = "A proposal for discussion"

# --- paragraph 1: "Initial Target milestone " -> "Need to find a balance between"
#     then insert two new level-1 paragraphs after it.
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "TEMP_PLACEHOLDER_1"
$para1b = $tr.Paragraphs(1, 1)
$para1b.Text = "Need to find a balance between"

$para1c = $tr.Paragraphs(1, 1)
[void]$para1c.InsertAfter("`rA minimal set of core terms use widely across all NMOP documents`rA maximal set of terms that might be used in any one NMOP document ")

$newPara2 = $tr.Paragraphs(2, 1)
$newPara2.IndentLevel = 2
$newPara3 = $tr.Paragraphs(3, 1)
$newPara3.IndentLevel = 2

$newPara2 = $tr.Paragraphs(2, 1)
$p2start = $newPara2.Start
$plain2 = "A minimal set of "
$bi2 = "core terms use widely across all"
$sub2 = $tr.Characters($p2start + $plain2.Length, $bi2.Length)
$sub2.Font.Bold = 1
$sub2.Font.Italic = 1

$newPara3 = $tr.Paragraphs(3, 1)
$p3start = $newPara3.Start
$plain3 = "A maximal set of terms that "
$bi3 = "might be used in any one "
$sub3 = $tr.Characters($p3start + $plain3.Length, $bi3.Length)
$sub3.Font.Bold = 1
$sub3.Font.Italic = 1

# ---------------------------------------------------------------------
# Reposition/resize the content placeholder (new explicit xfrm).
# Target EMU: off(838200,1581150) ext(10515600,4911725) -> points (/12700)
# ---------------------------------------------------------------------
$bodyShp.Left = 66.0
$bodyShp.Top = 124.5
$bodyShp.Width = 828.0
$bodyShp.Height = 386.75

# ---------------------------------------------------------------------
# Remove the "Tableau 5" table graphicFrame from the slide.
# ---------------------------------------------------------------------
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Tableau 5") {
        $shp.Delete()
    }
}
